$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.471.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +5.41%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.623.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +5.44%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'590.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.89%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'191.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.644"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.617.83"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +5.44%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.00%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.177"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.97%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +3.24%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'58.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +3.77%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +3.86%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +5.20%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.205.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +5.39%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'19.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +5.59%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.621.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +5.15%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'70.452.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +5.62%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'12.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.81%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.25%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +4.51%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'488.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.02%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'19.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +14.69%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'5.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +3.78%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +1.69%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'90.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.44%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +6.37%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'11.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.34%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'9.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +6.09%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'33.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +5.61%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +8.70%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'635.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +7.52%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +5.21%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +7.48%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'66.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.85%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'38.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +6.51%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +7.20%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +6.16%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.01%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.67%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'3.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.01%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.308.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +3.64%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'ThetaToken"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'3.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +6.30%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'Fetch.AI"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'2.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +11.07%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +5.19%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'3.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.90%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +2.89%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'9.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +4.51%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.61%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'3.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.62%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'Monero"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'142.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.67%  "
$ws.Range("E51").Style = "Normal"
